$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: extra hours logged against the weapon model / skin work ---
$ws.Range("C33").Formula = "=(1/60)*(21+8+10)"
$ws.Range("D33").Formula = "=(1/60)*(6)"
$ws.Range("E33").Formula = "=(1/60)*(9+13+21+21+21+21+11)"

# --- Row 34: C34 becomes an explicit (no-longer-shared) formula ---
$ws.Range("C34").Formula = "=(1/60)*(0)"

# --- New "BASED ON LAST DAY" metric (label in I8, value in I9) ---
$ws.Range("I8").Value = "BASED ON LAST DAY"
$ws.Range("I8").Font.Bold = $true
$ws.Range("I8").NumberFormat = "0.00"

$ws.Range("I9").Formula = "=I3/F33"

# --- Selection moves to D34 ---
$ws.Range("D34").Select() | Out-Null
